# Update interface definition calls in the "SOAPUI Tests" worksheet.
#
# The createOrder / confirmOrder sample SOAP messages are reshuffled:
#   - createOrder no longer sends shoppingCartInfo / shippingInfo inline
#     (those move down into the confirmOrder sample instead) and its
#     response now just returns a bare orderId.
#   - confirmOrder's purchaseOrder no longer carries a real accountId /
#     status, and gains the shoppingCartInfo / shippingInfo payloads.
#
# Also tidies a couple of long JSON cells (wrap + taller rows) and
# restores the worksheet's last saved scroll position/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SOAPUI Tests")

# --- createOrder block (rows 54-57) -------------------------------------
# C54 (Response) keeps the createOrderResponse opening tag
$ws.Range("C54").Value = "      <createOrderResponse xmlns=""http://OrderProcessService.WebServices.K9.com"">"

# Row 55
$ws.Range("B55").Value = "         <ord:shoppingCartInfo>[{""accountName"":""mbp"",""cdid"":1,""quantity"":3},{""accountName"":""mbp"",""cdid"":2,""quantity"":2}]</ord:shoppingCartInfo>"
$ws.Range("C55").Value = "         <createOrderReturn>{""orderId"":1}</createOrderReturn>"

# Row 56
$ws.Range("B56").Value = "         <ord:shippingInfo>{""accountName"":""mbp"",""shippingCharge"":5.25,""taxes"":4.25,""totalCost"":50.32}</ord:shippingInfo>"
$ws.Range("C56").Value = "      </createOrderResponse>"

# Row 57
$ws.Range("B57").Value = "      </ord:createOrder>"

# --- confirmOrder block (rows 61-67) -------------------------------------
$ws.Range("A61").Value = "confirmOrder"
$ws.Range("B62").Value = "<soapenv:Body>"

# Row 63
$ws.Range("B63").Value = "      <ord:confirmOrder>"
$ws.Range("C63").Value = "      <confirmOrderResponse xmlns=""http://OrderProcessService.WebServices.K9.com"">"

# Row 64
$ws.Range("B64").Value = "         <ord:purchaseOrder>{""orderId"":1,""accountId"":0,""status"":"""",""shippingCharge"":5.2,""taxes"":6.2,""totalCost"":20.36}</ord:purchaseOrder>"
$ws.Range("C64").Value = "         <confirmOrderReturn>{""callStatus"":0}</confirmOrderReturn>"

# Row 65
$ws.Range("B65").Value = "         <ord:shippingInfo>{""accountName"":""mbp"",""shippingCharge"":5.25,""taxes"":4.25,""totalCost"":50.32}</ord:shippingInfo>"
$ws.Range("C65").Value = "      </confirmOrderResponse>"

# Row 66
$ws.Range("B66").Value = "         <ord:paymentInfo>{""creditCardHolderName"":""MBP"",""creditCardNumber"":""4538452625981254"",""expiryDate"":""12/2018"",""ccv"":235</ord:paymentInfo>"

# Row 67
$ws.Range("B67").Value = "      </ord:confirmOrder>"

# --- row formatting -------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 60
$ws.Range("C5").WrapText = $true

$ws.Rows.Item(28).RowHeight = 105
$ws.Range("C28").WrapText = $true

# --- saved view / selection ------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 48
$ws.Range("C70").Select()
